# Apply the edits described by the diff to the "energy market" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values on row 2
$ws.Range("E2").Value = 278
$ws.Range("G2").Value = 501.41300000000012

# Update cell values on row 3
$ws.Range("D3").Value = 30850
$ws.Range("E3").Value = 280

# Update the active selection shown in the sheet view to H2
$ws.Range("H2").Select()
